$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "a"
$ws.Range("A2").Value = "b"

$ws.Range("A2").Select()
